# Edit script: insert new publication rows and reorganize existing ones
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 10; existing rows 10-15 shift down to 13-18
$ws.Rows("10:12").Insert()

# --- New row 10: https://openalex.org/W3039074189 ---
$ws.Range("A10").Value = '''https://openalex.org/W3039074189'
$ws.Range("B10").Value = '''Pattern of Road Traffic Accident Patients Admitted In Government Medical College and Hospital Nagpur – A Cross Sectional Study'
$ws.Range("C10").Value = '''list(au_id = c("https://openalex.org/A5040899345", "https://openalex.org/A5090400787", "https://openalex.org/A5047964975", "https://openalex.org/A5021789119"), au_display_name = c("Prafulla Sune", "Sandeep Bhelkar", "Uday Narlawar", "Sanjeev Chaudhary"), au_orcid = c(NA, NA, NA, NA), author_position = c("first", "middle", "middle", "last"), au_affiliation_raw = c("GMC Nagpur, Nagpur", "GMC Nagpur, Nagpur", "GMC Nagpur, Nagpur", "GMC Nagpur, Nagpur"), institution_id = c("https://openalex.org/I168974072", 
"https://openalex.org/I168974072", "https://openalex.org/I168974072", "https://openalex.org/I168974072"), institution_display_name = c("Government Medical College", "Government Medical College", "Government Medical College", "Government Medical College"), institution_ror = c("https://ror.org/026b7da27", "https://ror.org/026b7da27", "https://ror.org/026b7da27", "https://ror.org/026b7da27"), institution_country_code = c("IN", "IN", "IN", "IN"), institution_type = c("education", "education", "education", 
"education"), institution_lineage = c("https://openalex.org/I168974072", "https://openalex.org/I168974072", "https://openalex.org/I168974072", "https://openalex.org/I168974072"))'
$ws.Range("D10").Value = '''Introduction- Injuries are increasingly recognized as a global public health epidemic. Around the world, almost 16,000 people die every-day from all types of injuries. Injuries represent 12% of the global burden of disease, the third most important cause of overall mortality and the main cause of death among 1-40-year age groups.&#x0D; Methodology- This hospital based cross-sectional study was conducted among the road traffic accident patients admitted in trauma care center of Government Medical College and Hospital Nagpur.&#x0D; Result- Most common type of injury was abrasion in 91.45% followed by laceration in 79.61% of study subjects. Other common injuries were contusion, fracture, internal hemorrhage, crush injury and dislocation. Majority of the study subjects i.e. 67.11% had head injury.&#x0D; Conclusion- Head injury was the most common injury found in the study. Abrasion and laceration was also more commonly found in study subjects. Head injury was found more in non-users of personal protective devices which was statistically significant.'
$ws.Range("E10").Value = '''2022-03-10'
$ws.Range("F10").Value = '''National journal of community medicine'
$ws.Range("G10").Value = '''https://openalex.org/S2764358632'
$ws.Range("H10").Value = '''N/A'
$ws.Range("I10").Value = '''0976-3325'
$ws.Range("J10").Value = '''https://doi.org/10.5455/njcm.20200324100821'
$ws.Range("K10").Value = '''https://njcmindia.com/index.php/file/article/download/314/188'
$ws.Range("L10").Value = '''cc-by-sa'
$ws.Range("M10").Value = '''publishedVersion'
$ws.Range("N10").Value = '''240'
$ws.Range("O10").Value = '''243'
$ws.Range("P10").Value = '''11'
$ws.Range("Q10").Value = '''06'
$ws.Range("R10").Value = '''TRUE'
$ws.Range("S10").Value = '''TRUE'
$ws.Range("T10").Value = '''hybrid'
$ws.Range("U10").Value = '''https://njcmindia.com/index.php/file/article/download/314/188'
$ws.Range("V10").Value = '''FALSE'
$ws.Range("W10").Value = '''en'
$ws.Range("X10").Value = '''NA'
$ws.Range("Y10").Value = '''0'
$ws.Range("Z10").Value = '''2022'
$ws.Range("AA10").Value = '''https://api.openalex.org/works?filter=cites:W3039074189'
$ws.Range("AB10").Value = '''c(openalex = "https://openalex.org/W3039074189", doi = "https://doi.org/10.5455/njcm.20200324100821", mag = "3039074189")'
$ws.Range("AC10").Value = '''https://doi.org/10.5455/njcm.20200324100821'
$ws.Range("AD10").Value = '''article'
$ws.Range("AE10").Value = '''c("https://openalex.org/W2058268992", "https://openalex.org/W2103266699", "https://openalex.org/W2110206069", "https://openalex.org/W2132165965", "https://openalex.org/W2314697283", "https://openalex.org/W2323699862", "https://openalex.org/W2883058638")'
$ws.Range("AF10").Value = '''c("https://openalex.org/W2791959495", "https://openalex.org/W2070540847", "https://openalex.org/W2275554183", "https://openalex.org/W2278806788", "https://openalex.org/W3035266956", "https://openalex.org/W1481206241", "https://openalex.org/W4242352147", "https://openalex.org/W2170927500", "https://openalex.org/W1967302349", "https://openalex.org/W4281555576")'
$ws.Range("AG10").Value = '''FALSE'
$ws.Range("AH10").Value = '''FALSE'

# --- New row 11: https://openalex.org/W4210649205 ---
$ws.Range("A11").Value = '''https://openalex.org/W4210649205'
$ws.Range("B11").Value = '''Awareness and Satisfaction About COVAXIN Vaccination Services at an Immunization Clinic in Nagpur: A Cross-Sectional Study'
$ws.Range("C11").Value = '''list(au_id = c("https://openalex.org/A5072958078", "https://openalex.org/A5021789119", "https://openalex.org/A5031555037", "https://openalex.org/A5047964975", "https://openalex.org/A5070154438", "https://openalex.org/A5029243093"), au_display_name = c("Ujwala Ukey", "Sanjeev Chaudhary", "Sarita Sharma", "Uday Narlawar", "Ravikant Singh", "Aditi J Dabir"), au_orcid = c(NA, NA, "https://orcid.org/0000-0002-9720-2244", NA, NA, NA), author_position = c("first", "middle", "middle", "middle", "middle", 
"last"), au_affiliation_raw = c("", "", "", "", "", ""), institution_id = c(NA, NA, NA, NA, NA, NA), institution_display_name = c(NA, NA, NA, NA, NA, NA), institution_ror = c(NA, NA, NA, NA, NA, NA), institution_country_code = c(NA, NA, NA, NA, NA, NA), institution_type = c(NA, NA, NA, NA, NA, NA), institution_lineage = c(NA, NA, NA, NA, NA, NA))'
$ws.Range("D11").Value = '''Introduction Initially, coronavirus disease 2019 (COVID-19) vaccination was started in India for the elderly above 60 years of age. Adults with any comorbidity have been gradually included in the vaccination drive. It is empirical to gain insight into the satisfaction of these beneficiaries with the vaccination as it may act as an influencing factor for receiving the vaccine. Materials and methods This was a descriptive cross-sectional study carried out at the COVID-19 vaccination clinic of the Government Medical College and Hospital, Nagpur, among individuals above 60 years of age and those from 45 to 60 years of age with comorbidity. The survey tool was a predesigned structured questionnaire that had close-ended questions on various aspects of awareness about the COVID-19 vaccines and their satisfaction with the immunization center. Interviews were conducted by two interviewers on each day. Data were analyzed using open software Epi Info (CDC, Atlanta, Georgia). The chi-square test was applied as a test of significance. Results A total of 290 subjects participated in the study. The majority had correct knowledge about COVID-19 vaccination and appropriate COVID-19 behavior after vaccination. Fever and body ache were known to most of the subjects as adverse effects following immunization. Social media was the most common source of knowledge. The majority of the subjects were satisfied with the services provided at the vaccination center, but there was no difference as per age, gender, or residential status of the subjects. Conclusion Despite mixed rumors about the COVID-19 vaccine, the majority of the study subjects were well satisfied with the vaccination. They were apparently having fair awareness about the vaccine.'
$ws.Range("E11").Value = '''2022-01-06'
$ws.Range("F11").Value = '''Cureus'
$ws.Range("G11").Value = '''https://openalex.org/S2738950867'
$ws.Range("H11").Value = '''Cureus, Inc.'
$ws.Range("I11").Value = '''2168-8184'
$ws.Range("J11").Value = '''https://doi.org/10.7759/cureus.20983'
$ws.Range("K11").Value = '''https://www.cureus.com/articles/80272-awareness-and-satisfaction-about-covaxin-vaccination-services-at-an-immunization-clinic-in-nagpur-a-cross-sectional-study.pdf'
$ws.Range("L11").Value = '''N/A'
$ws.Range("M11").Value = '''publishedVersion'
$ws.Range("N11").Value = '''N/A'
$ws.Range("O11").Value = '''N/A'
$ws.Range("P11").Value = '''N/A'
$ws.Range("Q11").Value = '''N/A'
$ws.Range("R11").Value = '''TRUE'
$ws.Range("S11").Value = '''TRUE'
$ws.Range("T11").Value = '''gold'
$ws.Range("U11").Value = '''https://www.cureus.com/articles/80272-awareness-and-satisfaction-about-covaxin-vaccination-services-at-an-immunization-clinic-in-nagpur-a-cross-sectional-study.pdf'
$ws.Range("V11").Value = '''TRUE'
$ws.Range("W11").Value = '''en'
$ws.Range("X11").Value = '''NA'
$ws.Range("Y11").Value = '''0'
$ws.Range("Z11").Value = '''2022'
$ws.Range("AA11").Value = '''https://api.openalex.org/works?filter=cites:W4210649205'
$ws.Range("AB11").Value = '''c(openalex = "https://openalex.org/W4210649205", doi = "https://doi.org/10.7759/cureus.20983", pmid = "https://pubmed.ncbi.nlm.nih.gov/35154959")'
$ws.Range("AC11").Value = '''https://doi.org/10.7759/cureus.20983'
$ws.Range("AD11").Value = '''article'
$ws.Range("AE11").Value = '''c("https://openalex.org/W2766680409", "https://openalex.org/W3017185871", "https://openalex.org/W3046577186", "https://openalex.org/W3048424114", "https://openalex.org/W3093974194", "https://openalex.org/W3156243989", "https://openalex.org/W3158768260", "https://openalex.org/W3180633927", "https://openalex.org/W3181028174", "https://openalex.org/W3183958795", "https://openalex.org/W3187429287", "https://openalex.org/W3197228561", "https://openalex.org/W3208197226")'
$ws.Range("AF11").Value = '''c("https://openalex.org/W2981105526", "https://openalex.org/W2002551728", "https://openalex.org/W979319891", "https://openalex.org/W4386252105", "https://openalex.org/W2947320870", "https://openalex.org/W4367394835", "https://openalex.org/W2417808682", "https://openalex.org/W3135579318", "https://openalex.org/W3025685877", "https://openalex.org/W1929727402")'
$ws.Range("AG11").Value = '''FALSE'
$ws.Range("AH11").Value = '''FALSE'

# --- New row 12: https://openalex.org/W4220807774 ---
$ws.Range("A12").Value = '''https://openalex.org/W4220807774'
$ws.Range("B12").Value = '''Anaesthetic management of severe dextroscoliosis in a post-covid patient undergoing lumbar fixation and laminectomy'
$ws.Range("C12").Value = '''list(au_id = c("https://openalex.org/A5071817833", "https://openalex.org/A5091645233", "https://openalex.org/A5021789119", "https://openalex.org/A5041799937", "https://openalex.org/A5014814074"), au_display_name = c("Radhika Bajaj", "Amol Singam", "Sanjeev Chaudhary", "Rahul Chaudhary", "Parag Dongre"), au_orcid = c(NA, NA, NA, "https://orcid.org/0000-0002-3276-385X", NA), author_position = c("first", "middle", "middle", "middle", "last"), au_affiliation_raw = c("", "", "", "", ""), institution_id = c(NA, 
NA, NA, NA, NA), institution_display_name = c(NA, NA, NA, NA, NA), institution_ror = c(NA, NA, NA, NA, NA), institution_country_code = c(NA, NA, NA, NA, NA), institution_type = c(NA, NA, NA, NA, NA), institution_lineage = c(NA, NA, NA, NA, NA))'
$ws.Range("D12").Value = '''N/A'
$ws.Range("E12").Value = '''2022-03-11'
$ws.Range("F12").Value = '''Medical science'
$ws.Range("G12").Value = '''https://openalex.org/S4210211701'
$ws.Range("H12").Value = '''N/A'
$ws.Range("I12").Value = '''2321-7359'
$ws.Range("J12").Value = '''https://doi.org/10.54905/disssi/v26i121/ms97e2108'
$ws.Range("K12").Value = '''http://discoveryjournals.org/medicalscience/current_issue/v26/n121/ms97e2108.pdf#zoom=125'
$ws.Range("L12").Value = '''N/A'
$ws.Range("M12").Value = '''publishedVersion'
$ws.Range("N12").Value = '''1'
$ws.Range("O12").Value = '''1'
$ws.Range("P12").Value = '''26'
$ws.Range("Q12").Value = '''121'
$ws.Range("R12").Value = '''TRUE'
$ws.Range("S12").Value = '''TRUE'
$ws.Range("T12").Value = '''bronze'
$ws.Range("U12").Value = '''http://discoveryjournals.org/medicalscience/current_issue/v26/n121/ms97e2108.pdf#zoom=125'
$ws.Range("V12").Value = '''FALSE'
$ws.Range("W12").Value = '''en'
$ws.Range("X12").Value = '''NA'
$ws.Range("Y12").Value = '''0'
$ws.Range("Z12").Value = '''2022'
$ws.Range("AA12").Value = '''https://api.openalex.org/works?filter=cites:W4220807774'
$ws.Range("AB12").Value = '''c(openalex = "https://openalex.org/W4220807774", doi = "https://doi.org/10.54905/disssi/v26i121/ms97e2108")'
$ws.Range("AC12").Value = '''https://doi.org/10.54905/disssi/v26i121/ms97e2108'
$ws.Range("AD12").Value = '''article'
$ws.Range("AE12").Value = '''NA'
$ws.Range("AF12").Value = '''c("https://openalex.org/W4206669628", "https://openalex.org/W3198183218", "https://openalex.org/W3176864053", "https://openalex.org/W4382894326", "https://openalex.org/W3084498529", "https://openalex.org/W3020699490", "https://openalex.org/W4292098121", "https://openalex.org/W3036314732", "https://openalex.org/W4205317059", "https://openalex.org/W3009669391")'
$ws.Range("AG12").Value = '''FALSE'
$ws.Range("AH12").Value = '''FALSE'

# --- New row 19 (appended at end): https://openalex.org/W4387719557 ---
$ws.Range("A19").Value = '''https://openalex.org/W4387719557'
$ws.Range("B19").Value = '''Pattern of recurrent pediatric headache: A cohort of 100 children'
$ws.Range("C19").Value = '''list(au_id = c("https://openalex.org/A5010023158", "https://openalex.org/A5021789119", "https://openalex.org/A5029368727", "https://openalex.org/A5036686823", "https://openalex.org/A5051902511", "https://openalex.org/A5079035504"), au_display_name = c("Neerja Bhardwaj", "Sanjeev Chaudhary", "Amit Bhardwaj", "Neeraj Gupta", "VykuntarajuK Gowda", "AshwinVivek Sardesai"), au_orcid = c("https://orcid.org/0000-0001-5550-6277", NA, NA, "https://orcid.org/0000-0002-7131-4985", NA, NA), author_position = c("first", 
"middle", "middle", "middle", "middle", "last"), au_affiliation_raw = c("Department of Pediatrics, Dr. Rajendra Prasad Govt. Medical College Kangra at Tanda, Himachal Pradesh, India", "Department of Pediatrics, Dr. Rajendra Prasad Govt. Medical College Kangra at Tanda, Himachal Pradesh, India", "Department of Neurology, Dr. Rajendra Prasad Govt. Medical College Kangra at Tanda, Himachal Pradesh, India", "Department of Neonatology, AIIMS Jodhpur, Rajasthan, India", "Department of Pediatric Neurology, Indira Gandhi Institute of Child Health Bangalore, Karnataka, India", 
"Department of Pediatric Neurology, Indira Gandhi Institute of Child Health Bangalore, Karnataka, India"), institution_id = c("https://openalex.org/I4210148181", "https://openalex.org/I4210148181", "https://openalex.org/I4210148181", "https://openalex.org/I216021267", "https://openalex.org/I4210158200", "https://openalex.org/I4210158200"), institution_display_name = c("Dr. Rajendra Prasad Government Medical College", "Dr. Rajendra Prasad Government Medical College", "Dr. Rajendra Prasad Government Medical College", 
"All India Institute of Medical Sciences Jodhpur", "Indira Gandhi Institute of Child Health", "Indira Gandhi Institute of Child Health"), institution_ror = c("https://ror.org/04ce4rf90", "https://ror.org/04ce4rf90", "https://ror.org/04ce4rf90", "https://ror.org/05e15a779", "https://ror.org/04saq4y86", "https://ror.org/04saq4y86"), institution_country_code = c("IN", "IN", "IN", "IN", "IN", "IN"), institution_type = c("education", "education", "education", "education", "healthcare", "healthcare"), 
    institution_lineage = c("https://openalex.org/I4210148181", "https://openalex.org/I4210148181", "https://openalex.org/I4210148181", "https://openalex.org/I216021267, https://openalex.org/I4210148677", "https://openalex.org/I4210158200", "https://openalex.org/I4210158200"))'
$ws.Range("D19").Value = '''<b>Introduction:</b> This study was conducted to provide detailed information about clinical characteristics and short-term treatment outcome of childhood headache. <b>Materials and Methods:</b> This prospective observational study was done over a period of 15 months (January 2013 to March 2014) at a rural tertiary care center in North India. Detailed history, clinical examination, specialty review, and follow-up details were maintained on pretested structured proforma. Final diagnosis of headache type was made as per International Classification of Headache Disorders, 2nd edition. <b>Results:</b> Out of 100 (45 boys) children aged 8–18 years, 52% were diagnosed with migraine, 23% with tension-type headache, and 25% with secondary headache. Diffuse headache was the commonest (41%), and photophobia, phonophobia, and dizziness were the commonest symptoms in all headache subtypes. Sixty-five percentage of migraine headache were triggered by exertion. Lack of sleep and anxiety were triggers in most (65%) of tension type headache. Of 52 migraine children, 21 were started on prophylaxis for migraine and 14 of them reported significant improvement. Secondary causes for headache were found in 25% of children and half of them were having refractive errors. <b>Conclusions:</b> Results of the study show migraine being the commonest type of headache in children followed by secondary headache. This study also highlights the need for long-term follow-up of childhood headache.'
$ws.Range("E19").Value = '''2022-01-01'
$ws.Range("F19").Value = '''Journal of Pediatric Neurosciences'
$ws.Range("G19").Value = '''https://openalex.org/S195430055'
$ws.Range("H19").Value = '''Medknow'
$ws.Range("I19").Value = '''1817-1745'
$ws.Range("J19").Value = '''https://doi.org/10.4103/jpn.jpn_142_20'
$ws.Range("K19").Value = '''N/A'
$ws.Range("L19").Value = '''N/A'
$ws.Range("M19").Value = '''N/A'
$ws.Range("N19").Value = '''23'
$ws.Range("O19").Value = '''23'
$ws.Range("P19").Value = '''17'
$ws.Range("Q19").Value = '''1'
$ws.Range("R19").Value = '''FALSE'
$ws.Range("S19").Value = '''FALSE'
$ws.Range("T19").Value = '''closed'
$ws.Range("U19").Value = '''N/A'
$ws.Range("V19").Value = '''FALSE'
$ws.Range("W19").Value = '''en'
$ws.Range("X19").Value = '''NA'
$ws.Range("Y19").Value = '''0'
$ws.Range("Z19").Value = '''2022'
$ws.Range("AA19").Value = '''https://api.openalex.org/works?filter=cites:W4387719557'
$ws.Range("AB19").Value = '''c(openalex = "https://openalex.org/W4387719557", doi = "https://doi.org/10.4103/jpn.jpn_142_20")'
$ws.Range("AC19").Value = '''https://doi.org/10.4103/jpn.jpn_142_20'
$ws.Range("AD19").Value = '''article'
$ws.Range("AE19").Value = '''c("https://openalex.org/W1987154055", "https://openalex.org/W1987449740", "https://openalex.org/W1989499763", "https://openalex.org/W1989537855", "https://openalex.org/W1996136485", "https://openalex.org/W2007193001", "https://openalex.org/W2021662083", "https://openalex.org/W2062616749", "https://openalex.org/W2074256173", "https://openalex.org/W2076569501", "https://openalex.org/W2078144872", "https://openalex.org/W2084927512", "https://openalex.org/W2093370271", "https://openalex.org/W2106518663", 
"https://openalex.org/W2109212419", "https://openalex.org/W2132689324", "https://openalex.org/W2135718058", "https://openalex.org/W2152007006", "https://openalex.org/W2769724041", "https://openalex.org/W2800323269", "https://openalex.org/W2806669616", "https://openalex.org/W2885199469", "https://openalex.org/W2900091888")'
$ws.Range("AF19").Value = '''c("https://openalex.org/W1966605563", "https://openalex.org/W4207056969", "https://openalex.org/W2020802408", "https://openalex.org/W1978753422", "https://openalex.org/W181586006", "https://openalex.org/W4313043634", "https://openalex.org/W2029248794", "https://openalex.org/W1992468542", "https://openalex.org/W2329765567", "https://openalex.org/W2099842030")'
$ws.Range("AG19").Value = '''FALSE'
$ws.Range("AH19").Value = '''FALSE'

# Reset auto-calculated row heights back to default (matches source formatting)
$ws.Rows("10:12").AutoFit()
$ws.Rows("19:19").AutoFit()

